$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "66.533.15"
Set-TextValue "E2" "  -1.10%  "
Set-TextValue "D3" "3.451.74"
Set-TextValue "E3" "  -0.85%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "580.00"
Set-TextValue "E5" "  -2.37%  "
Set-TextValue "D6" "175.85"
Set-TextValue "E6" "  -1.53%  "
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "E8" "  +0.15%  "
Set-TextValue "D9" "3.449.55"
Set-TextValue "E9" "  -1.04%  "
Set-TextValue "E10" "  -2.20%  "
Set-TextValue "D11" "6.86"
Set-TextValue "E11" "  -3.20%  "
Set-TextValue "E12" "  -3.74%  "
Set-TextValue "D13" "4.044.82"
Set-TextValue "E13" "  -0.97%  "
Set-TextValue "D14" "30.44"
Set-TextValue "E14" "  -4.76%  "
Set-TextValue "E15" "  -3.37%  "
Set-TextValue "D16" "66.507.34"
Set-TextValue "E16" "  -1.25%  "
Set-TextValue "E17" "  -2.64%  "
Set-TextValue "D18" "3.450.27"
Set-TextValue "E18" "  -0.89%  "
Set-TextValue "E19" "  -3.98%  "
Set-TextValue "D20" "13.85"
Set-TextValue "E20" "  -3.06%  "
Set-TextValue "D21" "375.97"
Set-TextValue "E21" "  -3.29%  "
Set-TextValue "E22" "  -3.02%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.13%  "
Set-TextValue "E24" "  -0.22%  "
Set-TextValue "D25" "71.13"
Set-TextValue "E26" "  -1.52%  "
Set-TextValue "E27" "  -2.87%  "
Set-TextValue "D28" "9.78"
Set-TextValue "D29" "0.171"
Set-TextValue "E29" "  -2.13%  "
Set-TextValue "E30" "  +0.04%  "
Set-TextValue "E31" "  -5.14%  "
Set-TextValue "D32" "23.97"
Set-TextValue "E32" "  +1.94%  "
Set-TextValue "E33" "  -3.72%  "
Set-TextValue "E34" "  -5.93%  "
Set-TextValue "E36" "  -4.54%  "
Set-TextValue "E37" "  -5.00%  "
Set-TextValue "D38" "159.21"
Set-TextValue "E38" "  -2.89%  "
Set-TextValue "E39" "  +0.42%  "
Set-TextValue "E40" "  +4.32%  "
Set-TextValue "E41" "  -4.64%  "
Set-TextValue "D42" "2.63"
Set-TextValue "E42" "  -3.33%  "
Set-TextValue "D43" "6.50"
Set-TextValue "E44" "  -4.16%  "
Set-TextValue "D45" "2.682.72"
Set-TextValue "E45" "  -5.68%  "
Set-TextValue "E46" "  -3.83%  "
Set-TextValue "D47" "25.19"
Set-TextValue "E47" "  -6.34%  "
Set-TextValue "D48" "40.16"
Set-TextValue "E48" "  -3.40%  "
Set-TextValue "D49" "0.0295"
Set-TextValue "E49" "  -1.43%  "
Set-TextValue "D50" "319.75"
Set-TextValue "E50" "  -5.15%  "
Set-TextValue "E51" "  -4.19%  "
